# Switched LED footprint and fixed row traces
#  - Switched LED footprints from 0402 to 0603 to reduce assembly cost
#  - Added back row traces that were accidentally removed in an earlier commit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Added back row traces (C3/C4 caps folded into the C1/C2 group; new C5 row) ---
$ws.Range("B2").Value = "C1 C2 C3 C4"
$ws.Range("B3").Value = "C5"
$ws.Range("D3").Value = "C23733"

# --- Switched LED footprint from 0402 to 0603, updated JLCPCB part ---
$ws.Range("C4").Value = 603
$ws.Range("D4").Value = "C2286"

# --- Updated JLCPCB part numbers (resistor rows) ---
$ws.Range("D6").Value = "C25744"
$ws.Range("D7").Value = "C11702"
$ws.Range("D8").Value = "C25091"

# --- View / formatting touch-ups ---
$ws.Columns("A").ColumnWidth = 15.83
$excel.ActiveWindow.Zoom = 130
[void]$ws.Range("E4").Select()
